# ADI-scrum-BurndownChart.xlsx edit
# - Fill in column N ("actual burndown" running series for the latest day)
#   for rows 5-26 with the day's reported numbers, and total them in N29.
# - Update the sheet view (scroll position / active selection) left by the
#   user after finishing data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADI-burndown")
$ws.Activate()

# Per-task "remaining work" entries added in column N (mirrors the other
# day columns, e.g. column M) for rows 5 through 26.
$nValues = @{
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 30
    25 = 15
    26 = 10
}

foreach ($row in $nValues.Keys) {
    $ws.Range("N$row").Value = $nValues[$row]
}

# Column total (mirrors the SUM formula already present for the other
# day columns, e.g. M29).
$ws.Range("N29").Formula = "=SUM(N5:N27)"

# Restore the sheet view state (scroll position + active selection) as left
# by the user.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P29").Select()
